$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.143.77'
$ws.Range('E2').Value = '  -0.81%  '

$ws.Range('D3').Value = '2.026.01'
$ws.Range('E3').Value = '  -1.96%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.26'
$ws.Range('E5').Value = '  -1.98%  '

$ws.Range('E6').Value = '  -2.09%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '55.18'
$ws.Range('E8').Value = '  -4.46%  '

$ws.Range('E9').Value = '  -2.82%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0788'
$ws.Range('E10').Value = '  +0.59%  '

$ws.Range('E11').Value = '  -4.90%  '

$ws.Range('D12').Value = '2.323.29'
$ws.Range('E12').Value = '  -2.01%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.23'
$ws.Range('E13').Value = '  -4.00%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.35'
$ws.Range('E14').Value = '  -3.50%  '

$ws.Range('E15').Value = '  -2.65%  '

$ws.Range('E16').Value = '  -3.07%  '

$ws.Range('D17').Value = '2.017.01'
$ws.Range('E17').Value = '  -2.48%  '

$ws.Range('D18').Value = '37.149.40'
$ws.Range('E18').Value = '  -0.71%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.35'
$ws.Range('E19').Value = '  +3.48%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.14'
$ws.Range('E20').Value = '  -1.48%  '

$ws.Range('E21').Value = '  -1.22%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.16'
$ws.Range('E22').Value = '  -1.45%  '

$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('E24').Value = '  +1.75%  '

$ws.Range('E25').Value = '  -5.94%  '

$ws.Range('E26').Value = '  -6.58%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.76'
$ws.Range('E27').Value = '  -2.12%  '

$ws.Range('E28').Value = '  -3.38%  '

$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.35'
$ws.Range('E29').Value = '  -1.77%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '18.75'
$ws.Range('E30').Value = '  -2.77%  '

$ws.Range('E31').Value = '  -3.02%  '

$ws.Range('E32').Value = '  -1.57%  '

$ws.Range('E33').Value = '  -2.22%  '

$ws.Range('E34').Value = '  -4.77%  '

$ws.Range('E36').Value = '  +1.28%  '

$ws.Range('E37').Value = '  +0.23%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.16'
$ws.Range('E38').Value = '  -4.36%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.46'
$ws.Range('E39').Value = '  +2.87%  '

$ws.Range('E40').Value = '  -4.51%  '

$ws.Range('D41').Value = '1.477.00'
$ws.Range('E41').Value = '  -0.92%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '96.17'
$ws.Range('E42').Value = '  -2.00%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0919'
$ws.Range('E43').Value = '  -3.75%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.35'
$ws.Range('E44').Value = '  -3.30%  '

$ws.Range('E45').Value = '  -5.08%  '

$ws.Range('E46').Value = '  -4.65%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.28'
$ws.Range('E47').Value = '  +0.53%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.02'
$ws.Range('E48').Value = '  -2.05%  '

$ws.Range('E49').Value = '  -1.13%  '

$ws.Range('D50').Value = '2.212.18'
$ws.Range('E50').Value = '  -1.95%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.61'
$ws.Range('E51').Value = '  -10.84%  '
